# Start emulator and complete clock module
#
# Adds the "Upper MAR" / "Lower MAR" bit-weight label block (rows 24-27),
# mirroring the existing "Upper/Lower Temp 1", "Temp 2" and "ALU Result"
# blocks above it, and nudges a couple of cosmetic sheet properties
# (selection cell, column width for A:J).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Minor sheetView / column tweaks
# ---------------------------------------------------------------------
$ws.Range("U25").Select()
$ws.Columns("A:J").ColumnWidth = 1.31

# ---------------------------------------------------------------------
# 2. New title row (24): "Upper MAR" / "Lower MAR" headers.
#    Merge first, then copy the formatting from the row-19 title row
#    ("Upper ALU Result" / "Lower ALU Result") so borders/font/
#    alignment match exactly (same style ids as the existing blocks),
#    then set the text.
# ---------------------------------------------------------------------
$ws.Range("A24:J24").Merge()
$ws.Range("L24:U24").Merge()

$ws.Range("A19:J19").Copy()
$ws.Range("A24:J24").PasteSpecial(-4122) | Out-Null
$ws.Range("L19:U19").Copy()
$ws.Range("L24:U24").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("A24").Value = "Upper MAR"
$ws.Range("L24").Value = "Lower MAR"

# ---------------------------------------------------------------------
# 3. New data rows (25-27): bit-weight label table. Merge each column
#    across the 3 rows first (matching the Temp1/Temp2/ALU Result
#    blocks), then copy formatting from the row-20:22 data block, then
#    fill in the labels.
# ---------------------------------------------------------------------
$cols = @("A","B","C","D","E","F","G","H","I","J","L","M","N","O","P","Q","R","S","T","U")
foreach ($col in $cols) {
    $addr = $col + "25:" + $col + "27"
    $ws.Range($addr).Merge()
}

$ws.Range("A20:J22").Copy()
$ws.Range("A25:J27").PasteSpecial(-4122) | Out-Null
$ws.Range("L20:U22").Copy()
$ws.Range("L25:U27").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Upper MAR byte (A25:J25): 32K 16K 8K 4K 2K 1K 512 256 .. IN
$ws.Range("A25").Value = "3`n2`nK"
$ws.Range("B25").Value = "1`n6`nK"
$ws.Range("C25").Value = "8`nK"
$ws.Range("D25").Value = "4`nK"
$ws.Range("E25").Value = "2`nK"
$ws.Range("F25").Value = "1`nK"
$ws.Range("G25").Value = "5`n1`n2"
$ws.Range("H25").Value = "2`n5`n6"
$ws.Range("I25").Value = "I`nN"

# Lower MAR byte (L25:U25): 128 64 32 16 8 4 2 1 .. IN
$ws.Range("L25").Value = "1`n2`n8"
$ws.Range("M25").Value = "6`n4"
$ws.Range("N25").Value = "3`n2"
$ws.Range("O25").Value = "1`n6"
$ws.Range("P25").Value = 8
$ws.Range("Q25").Value = 4
$ws.Range("R25").Value = 2
$ws.Range("S25").Value = 1
$ws.Range("T25").Value = "I`nN"
